# Add a new "november-2025" worksheet at the end of the workbook, mirroring
# the structure of the existing monthly sheets (single cell A1 with a
# shared-string tax-revenue summary line).

$wb = $excel.ActiveWorkbook

# Insert the new worksheet after the current last sheet so it lands at the
# end of the tab order (matching sheetId 23 / rId23 placement in the diff).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "november-2025"

$ws.Range("A1").Value = ": tax revenue                                               96,968           99,100            -2,132            -2.2%"
